$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "A1" = "예윤 병국"; "B1" = "태훈 서연"; "C1" = "예윤 병국"; "D1" = "서연 한솔"; "E1" = "태훈 재현";
    "A2" = "예윤 유진"; "B2" = "태훈 서연"; "C2" = "예윤 희지"; "D2" = "현빈 태훈"; "E2" = "태훈 재현";
    "A3" = "병국 유진"; "B3" = "태훈 현빈"; "C3" = "희지 유진"; "D3" = "현빈 태훈"; "E3" = "재현 병국";
    "A4" = "유진 서연"; "B4" = "태훈 한솔"; "C4" = "유진 서연"; "D4" = "준범";     "E4" = "병국 희지";
    "A5" = "서연 한솔"; "B5" = "태훈 예윤"; "C5" = "서연 한솔"; "D5" = "준범 유진"; "E5" = "병국";
    "A6" = "서연 한솔"; "B6" = "예윤 현빈"; "C6" = "서연 한솔"; "D6" = "준범 희지"; "E6" = "병국 혜지";
    "A7" = "한솔 희지"; "B7" = "현빈 예윤"; "C7" = "준범 유진"; "D7" = "희지 현빈"; "E7" = "혜지";
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
